$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New pattern names to append (shared strings indices 20-30)
$newPatterns = @(
    "DomainExpertiseInRoles",
    "FeatureAssignment",
    "FewRoles",
    "GenericsAndSpecifics",
    "HierarchyOfFactories",
    "LockEmUpTogether",
    "LooseInterfaces",
    "OrganizationFollowsMarket",
    "ArchitectureTeam",
    "CodeOwnership",
    "DistributeWorkEvenly"
)

# Values used across all data rows (row 2-8 pattern repeats)
$rowValues = @(1.01, 1.01, 1.01, 1.01, 1.01, 0.89, 0.89, 0.89, 1.01, 1.01, 1.01)

$startRow = 9
for ($i = 0; $i -lt $newPatterns.Length; $i++) {
    $r = $startRow + $i

    # Duplicate the row immediately above (same formatting: borders,
    # alignment, fonts) into the new row so styles carry over faithfully.
    $ws.Rows.Item($r - 1).Copy()
    $ws.Rows.Item($r).Insert(-4121)
    # The insert step can lose the thin border on column A; restore it so
    # it resolves back to the same cell style used by the other rows.
    $ws.Range("A" + $r).Borders.LineStyle = 1

    $ws.Cells.Item($r, 1).Value = $newPatterns[$i]

    for ($c = 2; $c -le 12; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 2]
    }

    $ws.Cells.Item($r, 13).Value = ($r - 1)
}

$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("M19").Select()
